# Auto-generated Excel COM-interop script
# Updates the cryptos worksheet cell values per the commit diff.
#
# Cells whose new text looks like a plain number (e.g. "235.71") need a
# little care: assigning such a string straight to Range.Value lets Excel
# auto-type it into a real number (and changing a cells NumberFormat to
# force text first would bump its style index, which the source file never
# had). Instead we stage the text on a scratch cell formatted as Text, copy
# it, and Paste Special (values only) onto the destination cell -- that
# carries over the literal text without touching the destinations style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$scratch = $ws.Range("ZZ1")
$scratch.NumberFormat = "@"

$ws.Range("D2").Value = "25.611.53"
$ws.Range("E2").Value = "  -3.33%  "
$ws.Range("D3").Value = "1.735.92"
$ws.Range("E3").Value = "  -5.71%  "
$scratch.Value = "235.71"
$scratch.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  -10.03%  "
$ws.Range("E6").Value = "  +0.01%  "
$scratch.Value = "0.4887"
$scratch.Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").Value = "  -7.54%  "
$ws.Range("E8").Value = "  -8.33%  "
$scratch.Value = "0.2555"
$scratch.Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = "  -16.98%  "
$scratch.Value = "0.06089"
$scratch.Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = "  -11.59%  "
$ws.Range("D11").Value = "1.736.12"
$ws.Range("E11").Value = "  -5.73%  "
$scratch.Value = "0.06845"
$scratch.Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = "  -12.48%  "
$scratch.Value = "14.76"
$scratch.Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = "  -19.93%  "
$scratch.Value = "4.424"
$scratch.Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = "  -11.99%  "
$scratch.Value = "75.37"
$scratch.Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = "  -15.87%  "
$scratch.Value = "0.5623"
$scratch.Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Value = "  -25.79%  "
$ws.Range("E17").Value = "  +0.03%  "
$scratch.Value = "1.002"
$scratch.Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("D19").Value = "25.640.53"
$ws.Range("E19").Value = "  -3.33%  "
$scratch.Value = "11.42"
$scratch.Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = "  -18.39%  "
$scratch.Value = "0.000006562"
$scratch.Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = "  -17.39%  "
$ws.Range("D22").Value = "1.956.20"
$ws.Range("E22").Value = "  -5.92%  "
$scratch.Value = "4.019"
$scratch.Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = "  -13.03%  "
$scratch.Value = "7.896"
$scratch.Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = "  -15.11%  "
$scratch.Value = "4.980"
$scratch.Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = "  -16.94%  "
$scratch.Value = "136.64"
$scratch.Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = "  -3.92%  "
$scratch.Value = "1.490"
$scratch.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = "  -11.81%  "
$scratch.Value = "1.812"
$scratch.Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = "  -17.19%  "
$scratch.Value = "14.61"
$scratch.Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = "  -14.01%  "
$scratch.Value = "101.23"
$scratch.Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = "  -8.84%  "
$scratch.Value = "0.07983"
$scratch.Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = "  -9.10%  "
$scratch.Value = "3.667"
$scratch.Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = "  -14.07%  "
$scratch.Value = "3.379"
$scratch.Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = "  -17.27%  "
$scratch.Value = "0.04394"
$scratch.Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = "  -8.96%  "
$scratch.Value = "1.000"
$scratch.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = "  +0.00%  "
$scratch.Value = "2.613"
$scratch.Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = "  -10.91%  "
$scratch.Value = "0.9483"
$scratch.Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = "  -16.21%  "
$scratch.Value = "0.5874"
$scratch.Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = "  -19.49%  "
$scratch.Value = "2.624"
$scratch.Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = "  -15.35%  "
$ws.Range("E40").Value = "  +0.14%  "
$scratch.Value = "103.29"
$scratch.Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = "  -4.38%  "
$scratch.Value = "0.01496"
$scratch.Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = "  -13.00%  "
$scratch.Value = "1.856"
$scratch.Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = "  -19.62%  "
$scratch.Value = "5.113"
$scratch.Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = "  -13.10%  "
$scratch.Value = "0.3722"
$scratch.Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = "  -22.38%  "
$scratch.Value = "0.7216"
$scratch.Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = "  -20.01%  "
$scratch.Value = "0.05205"
$scratch.Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = "  -10.24%  "
$scratch.Value = "0.1082"
$scratch.Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = "  -12.66%  "
$scratch.Value = "29.79"
$scratch.Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = "  -14.66%  "
$ws.Range("E50").Value = "  -14.40%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$scratch.Value = "1.002"
$scratch.Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = "  +0.07%  "

$scratch.Clear()
$excel.CutCopyMode = 0
